$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row -> column value updates (only cells whose value changed per the diff)
# Columns: C (price), D (change %), I (dividend yield), J (%K), K (%D)

$updates = @{
    2  = @{ C=19780;  D=0.0302;  I=5.0599999999999996; J=70;  K=70 }
    3  = @{ C=98700;  D=0.0186;  I=6.59;               J=62;  K=62 }
    4  = @{ C=439500; D=0.0127;  I=4.32;               J=91;  K=91 }
    5  = @{ C=32000;  D=0.0273;  I=6.25;               J=50;  K=50 }
    6  = @{ C=30750;  D=-0.0653; I=3.9;                J=79;  K=79 }
    7  = @{ C=25150;  D=0.0307;  I=4.7699999999999996; J=91;  K=91 }
    8  = @{ C=10190;              I=5.05;               J=93;  K=93 }
    9  = @{ C=78000;  D=0.0223;  I=3.85;               J=65;  K=65 }
    10 = @{ C=210000; D=0.0319;  I=5.71;               J=49;  K=49 }
    11 = @{ C=125800; D=0.0203;  I=5.41;               J=96;  K=96 }
    12 = @{ C=20100;  D=0.0055;  I=4.7300000000000004; J=92;  K=92 }
    13 = @{ C=73400;  D=-0.0027; I=4.7699999999999996; J=89;  K=89 }
    14 = @{ C=57500;  D=0.0141;  I=6.16;               J=79;  K=79 }
    15 = @{ C=87300;  D=0.0151;  I=6.3;                J=94;  K=94 }
    16 = @{ C=18470;  D=0.0115;  I=5.77;               J=95;  K=95 }
    17 = @{ C=50500;  D=0.0223;  I=5.54;               J=94;  K=94 }
    18 = @{ C=20200;  D=0.0075;  I=6.09;               J=37;  K=37 }
    19 = @{ C=57300;  D=0.0287;  I=3.49;               J=100; K=100 }
    20 = @{ C=14380;  D=0.0056;  I=4.5199999999999996; J=73;  K=73 }
    21 = @{ C=127200; D=-0.0039; I=4.25;               J=97;  K=97 }
    22 = @{ C=41500;  D=-0.0036; I=3.51 }
    23 = @{ C=62200;  D=0.013;   I=3.47;               J=94;  K=94 }
    24 = @{ C=48750;  D=0.045;   I=5.54;               J=68;  K=68 }
    25 = @{ C=88900;  D=0.0301;  I=4.05;               J=97;  K=97 }
    26 = @{ C=110800; D=-0.0009 }
    27 = @{ C=12820;  D=0.0297;  I=5.07;               J=97;  K=97 }
    28 = @{ C=13000;  D=0.0342;  I=3.85;               J=98;  K=98 }
    29 = @{ C=20650;  D=0.0049;  I=4.82;               J=91;  K=91 }
}

foreach ($r in $updates.Keys) {
    $row = $updates[$r]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}

# Row 7: number format changes from 0% (no decimals) to 0.00% (two decimals)
$ws.Range("D7").NumberFormat = "0.00%"

# Row 30: D changes value to 0 and number format changes from 0.00% to 0%
$ws.Range("D30").Value = 0
$ws.Range("D30").NumberFormat = "0%"

# Update the active cell selection to F15 (as reflected in the sheetView)
$ws.Range("F15").Select()
